$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(68, 4).Value = 44447
$ws.Cells.Item(68, 12).Value = "Primera"
$ws.Cells.Item(68, 13).Value = 120
$ws.Cells.Item(68, 14).Value = 12000
$ws.Cells.Item(68, 15).Value = 12500
$ws.Cells.Item(68, 16).Value = 12250
$ws.Cells.Item(68, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(68, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(68, 19).Value = 681
$ws.Cells.Item(68, 20).Value = 18

$ws.Cells.Item(69, 4).Value = 44447
$ws.Cells.Item(69, 12).Value = "Segunda"
$ws.Cells.Item(69, 13).Value = 60
$ws.Cells.Item(69, 14).Value = 11000
$ws.Cells.Item(69, 15).Value = 11500
$ws.Cells.Item(69, 16).Value = 11250
$ws.Cells.Item(69, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(69, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(69, 19).Value = 625
$ws.Cells.Item(69, 20).Value = 18

$ws.Cells.Item(70, 4).Value = 44392
$ws.Cells.Item(70, 12).Value = "Primera"
$ws.Cells.Item(70, 13).Value = 120
$ws.Cells.Item(70, 14).Value = 10000
$ws.Cells.Item(70, 15).Value = 11000
$ws.Cells.Item(70, 16).Value = 10500
$ws.Cells.Item(70, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(70, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(70, 19).Value = 583
$ws.Cells.Item(70, 20).Value = 18

$ws.Cells.Item(71, 4).Value = 44392
$ws.Cells.Item(71, 12).Value = "Segunda"
$ws.Cells.Item(71, 13).Value = 80
$ws.Cells.Item(71, 14).Value = 9000
$ws.Cells.Item(71, 15).Value = 9000
$ws.Cells.Item(71, 16).Value = 9000
$ws.Cells.Item(71, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(71, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(71, 19).Value = 500
$ws.Cells.Item(71, 20).Value = 18

$ws.Cells.Item(72, 4).Value = 44362
$ws.Cells.Item(72, 12).Value = "Primera"
$ws.Cells.Item(72, 13).Value = 160
$ws.Cells.Item(72, 14).Value = 8000
$ws.Cells.Item(72, 15).Value = 8500
$ws.Cells.Item(72, 16).Value = 8250
$ws.Cells.Item(72, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(72, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(72, 19).Value = 825
$ws.Cells.Item(72, 20).Value = 10

$ws.Cells.Item(73, 4).Value = 44362
$ws.Cells.Item(73, 12).Value = "Segunda"
$ws.Cells.Item(73, 13).Value = 120
$ws.Cells.Item(73, 14).Value = 7000
$ws.Cells.Item(73, 15).Value = 7500
$ws.Cells.Item(73, 16).Value = 7250
$ws.Cells.Item(73, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(73, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(73, 19).Value = 725
$ws.Cells.Item(73, 20).Value = 10

$ws.Cells.Item(74, 4).Value = 44384
$ws.Cells.Item(74, 12).Value = "Primera"
$ws.Cells.Item(74, 13).Value = 120
$ws.Cells.Item(74, 14).Value = 10000
$ws.Cells.Item(74, 15).Value = 11000
$ws.Cells.Item(74, 16).Value = 10500
$ws.Cells.Item(74, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(74, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(74, 19).Value = 583
$ws.Cells.Item(74, 20).Value = 18

$ws.Cells.Item(75, 4).Value = 44384
$ws.Cells.Item(75, 12).Value = "Segunda"
$ws.Cells.Item(75, 13).Value = 80
$ws.Cells.Item(75, 14).Value = 9000
$ws.Cells.Item(75, 15).Value = 9000
$ws.Cells.Item(75, 16).Value = 9000
$ws.Cells.Item(75, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(75, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(75, 19).Value = 500
$ws.Cells.Item(75, 20).Value = 18

$ws.Cells.Item(76, 4).Value = 44365
$ws.Cells.Item(76, 12).Value = "Primera"
$ws.Cells.Item(76, 13).Value = 240
$ws.Cells.Item(76, 14).Value = 7000
$ws.Cells.Item(76, 15).Value = 8500
$ws.Cells.Item(76, 16).Value = 7750
$ws.Cells.Item(76, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(76, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(76, 19).Value = 775
$ws.Cells.Item(76, 20).Value = 10

$ws.Cells.Item(77, 4).Value = 44358
$ws.Cells.Item(77, 12).Value = "Primera"
$ws.Cells.Item(77, 13).Value = 120
$ws.Cells.Item(77, 14).Value = 10500
$ws.Cells.Item(77, 15).Value = 11000
$ws.Cells.Item(77, 16).Value = 10750
$ws.Cells.Item(77, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(77, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(77, 19).Value = 597
$ws.Cells.Item(77, 20).Value = 18

$ws.Cells.Item(78, 4).Value = 44358
$ws.Cells.Item(78, 12).Value = "Segunda"
$ws.Cells.Item(78, 13).Value = 120
$ws.Cells.Item(78, 14).Value = 8500
$ws.Cells.Item(78, 15).Value = 9000
$ws.Cells.Item(78, 16).Value = 8750
$ws.Cells.Item(78, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(78, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(78, 19).Value = 486
$ws.Cells.Item(78, 20).Value = 18

$ws.Cells.Item(79, 4).Value = 44433
$ws.Cells.Item(79, 12).Value = "Primera"
$ws.Cells.Item(79, 13).Value = 80
$ws.Cells.Item(79, 14).Value = 12000
$ws.Cells.Item(79, 15).Value = 12500
$ws.Cells.Item(79, 16).Value = 12250
$ws.Cells.Item(79, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(79, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(79, 19).Value = 681
$ws.Cells.Item(79, 20).Value = 18

$ws.Cells.Item(80, 4).Value = 44433
$ws.Cells.Item(80, 12).Value = "Segunda"
$ws.Cells.Item(80, 13).Value = 60
$ws.Cells.Item(80, 14).Value = 11000
$ws.Cells.Item(80, 15).Value = 11500
$ws.Cells.Item(80, 16).Value = 11250
$ws.Cells.Item(80, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(80, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(80, 19).Value = 625
$ws.Cells.Item(80, 20).Value = 18

$ws.Cells.Item(81, 4).Value = 44397
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 120
$ws.Cells.Item(81, 14).Value = 10000
$ws.Cells.Item(81, 15).Value = 11000
$ws.Cells.Item(81, 16).Value = 10500
$ws.Cells.Item(81, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(81, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(81, 19).Value = 583
$ws.Cells.Item(81, 20).Value = 18

$ws.Cells.Item(82, 4).Value = 44397
$ws.Cells.Item(82, 12).Value = "Segunda"
$ws.Cells.Item(82, 13).Value = 120
$ws.Cells.Item(82, 14).Value = 8000
$ws.Cells.Item(82, 15).Value = 9000
$ws.Cells.Item(82, 16).Value = 8500
$ws.Cells.Item(82, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(82, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(82, 19).Value = 472
$ws.Cells.Item(82, 20).Value = 18

$ws.Cells.Item(83, 4).Value = 44321
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 120
$ws.Cells.Item(83, 14).Value = 9500
$ws.Cells.Item(83, 15).Value = 10000
$ws.Cells.Item(83, 16).Value = 9750
$ws.Cells.Item(83, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(83, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(83, 19).Value = 975
$ws.Cells.Item(83, 20).Value = 10

$ws.Cells.Item(84, 4).Value = 44321
$ws.Cells.Item(84, 12).Value = "Segunda"
$ws.Cells.Item(84, 13).Value = 160
$ws.Cells.Item(84, 14).Value = 8000
$ws.Cells.Item(84, 15).Value = 8500
$ws.Cells.Item(84, 16).Value = 8250
$ws.Cells.Item(84, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(84, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(84, 19).Value = 825
$ws.Cells.Item(84, 20).Value = 10

$ws.Cells.Item(85, 4).Value = 44426
$ws.Cells.Item(85, 12).Value = "Primera"
$ws.Cells.Item(85, 13).Value = 120
$ws.Cells.Item(85, 14).Value = 12000
$ws.Cells.Item(85, 15).Value = 12500
$ws.Cells.Item(85, 16).Value = 12250
$ws.Cells.Item(85, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(85, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(85, 19).Value = 681
$ws.Cells.Item(85, 20).Value = 18

$ws.Cells.Item(86, 4).Value = 44426
$ws.Cells.Item(86, 12).Value = "Segunda"
$ws.Cells.Item(86, 13).Value = 60
$ws.Cells.Item(86, 14).Value = 11000
$ws.Cells.Item(86, 15).Value = 11500
$ws.Cells.Item(86, 16).Value = 11250
$ws.Cells.Item(86, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(86, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(86, 19).Value = 625
$ws.Cells.Item(86, 20).Value = 18

$ws.Cells.Item(87, 4).Value = 44421
$ws.Cells.Item(87, 12).Value = "Primera"
$ws.Cells.Item(87, 13).Value = 160
$ws.Cells.Item(87, 14).Value = 12000
$ws.Cells.Item(87, 15).Value = 12500
$ws.Cells.Item(87, 16).Value = 12250
$ws.Cells.Item(87, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(87, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(87, 19).Value = 681
$ws.Cells.Item(87, 20).Value = 18

$ws.Cells.Item(88, 4).Value = 44421
$ws.Cells.Item(88, 12).Value = "Segunda"
$ws.Cells.Item(88, 13).Value = 100
$ws.Cells.Item(88, 14).Value = 11000
$ws.Cells.Item(88, 15).Value = 11500
$ws.Cells.Item(88, 16).Value = 11250
$ws.Cells.Item(88, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(88, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(88, 19).Value = 625
$ws.Cells.Item(88, 20).Value = 18

$ws.Cells.Item(89, 4).Value = 44434
$ws.Cells.Item(89, 12).Value = "Primera"
$ws.Cells.Item(89, 13).Value = 120
$ws.Cells.Item(89, 14).Value = 12000
$ws.Cells.Item(89, 15).Value = 12500
$ws.Cells.Item(89, 16).Value = 12250
$ws.Cells.Item(89, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(89, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(89, 19).Value = 681
$ws.Cells.Item(89, 20).Value = 18

$ws.Cells.Item(90, 4).Value = 44434
$ws.Cells.Item(90, 12).Value = "Segunda"
$ws.Cells.Item(90, 13).Value = 100
$ws.Cells.Item(90, 14).Value = 11000
$ws.Cells.Item(90, 15).Value = 11500
$ws.Cells.Item(90, 16).Value = 11250
$ws.Cells.Item(90, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(90, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(90, 19).Value = 625
$ws.Cells.Item(90, 20).Value = 18

$ws.Cells.Item(91, 4).Value = 44420
$ws.Cells.Item(91, 12).Value = "Primera"
$ws.Cells.Item(91, 13).Value = 200
$ws.Cells.Item(91, 14).Value = 12000
$ws.Cells.Item(91, 15).Value = 12500
$ws.Cells.Item(91, 16).Value = 12250
$ws.Cells.Item(91, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(91, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(91, 19).Value = 681
$ws.Cells.Item(91, 20).Value = 18

$ws.Cells.Item(92, 4).Value = 44420
$ws.Cells.Item(92, 12).Value = "Segunda"
$ws.Cells.Item(92, 13).Value = 120
$ws.Cells.Item(92, 14).Value = 11000
$ws.Cells.Item(92, 15).Value = 11500
$ws.Cells.Item(92, 16).Value = 11250
$ws.Cells.Item(92, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(92, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(92, 19).Value = 625
$ws.Cells.Item(92, 20).Value = 18

$ws.Cells.Item(93, 4).Value = 44348
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 120
$ws.Cells.Item(93, 14).Value = 10000
$ws.Cells.Item(93, 15).Value = 11000
$ws.Cells.Item(93, 16).Value = 10500
$ws.Cells.Item(93, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(93, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(93, 19).Value = 583
$ws.Cells.Item(93, 20).Value = 18

$ws.Cells.Item(94, 4).Value = 44348
$ws.Cells.Item(94, 12).Value = "Segunda"
$ws.Cells.Item(94, 13).Value = 120
$ws.Cells.Item(94, 14).Value = 8000
$ws.Cells.Item(94, 15).Value = 9000
$ws.Cells.Item(94, 16).Value = 8500
$ws.Cells.Item(94, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(94, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(94, 19).Value = 472
$ws.Cells.Item(94, 20).Value = 18

$ws.Cells.Item(95, 4).Value = 44316
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 120
$ws.Cells.Item(95, 14).Value = 10000
$ws.Cells.Item(95, 15).Value = 11000
$ws.Cells.Item(95, 16).Value = 10500
$ws.Cells.Item(95, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(95, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(95, 19).Value = 1050
$ws.Cells.Item(95, 20).Value = 10

$ws.Cells.Item(96, 4).Value = 44316
$ws.Cells.Item(96, 12).Value = "Segunda"
$ws.Cells.Item(96, 13).Value = 120
$ws.Cells.Item(96, 14).Value = 8500
$ws.Cells.Item(96, 15).Value = 9000
$ws.Cells.Item(96, 16).Value = 8750
$ws.Cells.Item(96, 17).Value = "`$/caja 10 kilos"
$ws.Cells.Item(96, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(96, 19).Value = 875
$ws.Cells.Item(96, 20).Value = 10

$ws.Cells.Item(97, 4).Value = 44427
$ws.Cells.Item(97, 12).Value = "Primera"
$ws.Cells.Item(97, 13).Value = 160
$ws.Cells.Item(97, 14).Value = 12000
$ws.Cells.Item(97, 15).Value = 12500
$ws.Cells.Item(97, 16).Value = 12250
$ws.Cells.Item(97, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(97, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(97, 19).Value = 681
$ws.Cells.Item(97, 20).Value = 18

$ws.Cells.Item(98, 4).Value = 44427
$ws.Cells.Item(98, 12).Value = "Segunda"
$ws.Cells.Item(98, 13).Value = 120
$ws.Cells.Item(98, 14).Value = 11000
$ws.Cells.Item(98, 15).Value = 11500
$ws.Cells.Item(98, 16).Value = 11250
$ws.Cells.Item(98, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(98, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(98, 19).Value = 625
$ws.Cells.Item(98, 20).Value = 18

$ws.Cells.Item(99, 4).Value = 44414
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 300
$ws.Cells.Item(99, 14).Value = 12500
$ws.Cells.Item(99, 15).Value = 13000
$ws.Cells.Item(99, 16).Value = 12750
$ws.Cells.Item(99, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(99, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(99, 19).Value = 708
$ws.Cells.Item(99, 20).Value = 18

$ws.Cells.Item(100, 4).Value = 44414
$ws.Cells.Item(100, 12).Value = "Segunda"
$ws.Cells.Item(100, 13).Value = 120
$ws.Cells.Item(100, 14).Value = 11000
$ws.Cells.Item(100, 15).Value = 11500
$ws.Cells.Item(100, 16).Value = 11250
$ws.Cells.Item(100, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(100, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(100, 19).Value = 625
$ws.Cells.Item(100, 20).Value = 18

$ws.Cells.Item(101, 4).Value = 44417
$ws.Cells.Item(101, 12).Value = "Primera"
$ws.Cells.Item(101, 13).Value = 160
$ws.Cells.Item(101, 14).Value = 12500
$ws.Cells.Item(101, 15).Value = 13000
$ws.Cells.Item(101, 16).Value = 12750
$ws.Cells.Item(101, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(101, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(101, 19).Value = 708
$ws.Cells.Item(101, 20).Value = 18

$ws.Cells.Item(102, 4).Value = 44417
$ws.Cells.Item(102, 12).Value = "Segunda"
$ws.Cells.Item(102, 13).Value = 80
$ws.Cells.Item(102, 14).Value = 11000
$ws.Cells.Item(102, 15).Value = 11500
$ws.Cells.Item(102, 16).Value = 11250
$ws.Cells.Item(102, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(102, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(102, 19).Value = 625
$ws.Cells.Item(102, 20).Value = 18

$ws.Cells.Item(103, 4).Value = 44441
$ws.Cells.Item(103, 12).Value = "Primera"
$ws.Cells.Item(103, 13).Value = 120
$ws.Cells.Item(103, 14).Value = 12000
$ws.Cells.Item(103, 15).Value = 12500
$ws.Cells.Item(103, 16).Value = 12250
$ws.Cells.Item(103, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(103, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(103, 19).Value = 681
$ws.Cells.Item(103, 20).Value = 18

$ws.Cells.Item(104, 4).Value = 44441
$ws.Cells.Item(104, 12).Value = "Segunda"
$ws.Cells.Item(104, 13).Value = 60
$ws.Cells.Item(104, 14).Value = 11000
$ws.Cells.Item(104, 15).Value = 11500
$ws.Cells.Item(104, 16).Value = 11250
$ws.Cells.Item(104, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(104, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(104, 19).Value = 625
$ws.Cells.Item(104, 20).Value = 18

$ws.Cells.Item(105, 4).Value = 44432
$ws.Cells.Item(105, 12).Value = "Primera"
$ws.Cells.Item(105, 13).Value = 120
$ws.Cells.Item(105, 14).Value = 12000
$ws.Cells.Item(105, 15).Value = 12500
$ws.Cells.Item(105, 16).Value = 12250
$ws.Cells.Item(105, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(105, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(105, 19).Value = 681
$ws.Cells.Item(105, 20).Value = 18

$ws.Cells.Item(106, 4).Value = 44432
$ws.Cells.Item(106, 12).Value = "Segunda"
$ws.Cells.Item(106, 13).Value = 60
$ws.Cells.Item(106, 14).Value = 11000
$ws.Cells.Item(106, 15).Value = 11500
$ws.Cells.Item(106, 16).Value = 11250
$ws.Cells.Item(106, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(106, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(106, 19).Value = 625
$ws.Cells.Item(106, 20).Value = 18

$ws.Cells.Item(107, 4).Value = 44369
$ws.Cells.Item(107, 12).Value = "Primera"
$ws.Cells.Item(107, 13).Value = 120
$ws.Cells.Item(107, 14).Value = 12000
$ws.Cells.Item(107, 15).Value = 13000
$ws.Cells.Item(107, 16).Value = 12500
$ws.Cells.Item(107, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(107, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(107, 19).Value = 694
$ws.Cells.Item(107, 20).Value = 18

$ws.Cells.Item(108, 4).Value = 44369
$ws.Cells.Item(108, 12).Value = "Segunda"
$ws.Cells.Item(108, 13).Value = 120
$ws.Cells.Item(108, 14).Value = 10000
$ws.Cells.Item(108, 15).Value = 11000
$ws.Cells.Item(108, 16).Value = 10500
$ws.Cells.Item(108, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(108, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(108, 19).Value = 583
$ws.Cells.Item(108, 20).Value = 18

$ws.Cells.Item(109, 4).Value = 44379
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 120
$ws.Cells.Item(109, 14).Value = 9500
$ws.Cells.Item(109, 15).Value = 10000
$ws.Cells.Item(109, 16).Value = 9750
$ws.Cells.Item(109, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(109, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(109, 19).Value = 542
$ws.Cells.Item(109, 20).Value = 18

$ws.Cells.Item(110, 4).Value = 44379
$ws.Cells.Item(110, 12).Value = "Segunda"
$ws.Cells.Item(110, 13).Value = 120
$ws.Cells.Item(110, 14).Value = 8000
$ws.Cells.Item(110, 15).Value = 8500
$ws.Cells.Item(110, 16).Value = 8250
$ws.Cells.Item(110, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(110, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(110, 19).Value = 458
$ws.Cells.Item(110, 20).Value = 18

$ws.Cells.Item(111, 4).Value = 44315
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 160
$ws.Cells.Item(111, 14).Value = 8000
$ws.Cells.Item(111, 15).Value = 9000
$ws.Cells.Item(111, 16).Value = 8500
$ws.Cells.Item(111, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(111, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(111, 19).Value = 850
$ws.Cells.Item(111, 20).Value = 10

$ws.Cells.Item(112, 4).Value = 44391
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 120
$ws.Cells.Item(112, 14).Value = 10000
$ws.Cells.Item(112, 15).Value = 11000
$ws.Cells.Item(112, 16).Value = 10500
$ws.Cells.Item(112, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(112, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(112, 19).Value = 583
$ws.Cells.Item(112, 20).Value = 18

$ws.Cells.Item(113, 4).Value = 44446
$ws.Cells.Item(113, 12).Value = "Primera"
$ws.Cells.Item(113, 13).Value = 120
$ws.Cells.Item(113, 14).Value = 12000
$ws.Cells.Item(113, 15).Value = 12500
$ws.Cells.Item(113, 16).Value = 12250
$ws.Cells.Item(113, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(113, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(113, 19).Value = 681
$ws.Cells.Item(113, 20).Value = 18

$ws.Cells.Item(114, 4).Value = 44446
$ws.Cells.Item(114, 12).Value = "Segunda"
$ws.Cells.Item(114, 13).Value = 60
$ws.Cells.Item(114, 14).Value = 11000
$ws.Cells.Item(114, 15).Value = 11500
$ws.Cells.Item(114, 16).Value = 11250
$ws.Cells.Item(114, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(114, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(114, 19).Value = 625
$ws.Cells.Item(114, 20).Value = 18

$ws.Cells.Item(115, 4).Value = 44411
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 120
$ws.Cells.Item(115, 14).Value = 8000
$ws.Cells.Item(115, 15).Value = 8500
$ws.Cells.Item(115, 16).Value = 8250
$ws.Cells.Item(115, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(115, 18).Value = "Región del Maule"
$ws.Cells.Item(115, 19).Value = 825
$ws.Cells.Item(115, 20).Value = 10

$ws.Cells.Item(116, 4).Value = 44411
$ws.Cells.Item(116, 12).Value = "Segunda"
$ws.Cells.Item(116, 13).Value = 100
$ws.Cells.Item(116, 14).Value = 7000
$ws.Cells.Item(116, 15).Value = 7500
$ws.Cells.Item(116, 16).Value = 7250
$ws.Cells.Item(116, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(116, 18).Value = "Región del Maule"
$ws.Cells.Item(116, 19).Value = 725
$ws.Cells.Item(116, 20).Value = 10

$ws.Cells.Item(117, 4).Value = 44313
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 120
$ws.Cells.Item(117, 14).Value = 8000
$ws.Cells.Item(117, 15).Value = 9000
$ws.Cells.Item(117, 16).Value = 8500
$ws.Cells.Item(117, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(117, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(117, 19).Value = 850
$ws.Cells.Item(117, 20).Value = 10

$ws.Cells.Item(118, 4).Value = 44438
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 120
$ws.Cells.Item(118, 14).Value = 12000
$ws.Cells.Item(118, 15).Value = 12500
$ws.Cells.Item(118, 16).Value = 12250
$ws.Cells.Item(118, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(118, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(118, 19).Value = 681
$ws.Cells.Item(118, 20).Value = 18

$ws.Cells.Item(119, 4).Value = 44438
$ws.Cells.Item(119, 12).Value = "Segunda"
$ws.Cells.Item(119, 13).Value = 60
$ws.Cells.Item(119, 14).Value = 11000
$ws.Cells.Item(119, 15).Value = 11500
$ws.Cells.Item(119, 16).Value = 11250
$ws.Cells.Item(119, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(119, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(119, 19).Value = 625
$ws.Cells.Item(119, 20).Value = 18

$ws.Cells.Item(120, 4).Value = 44442
$ws.Cells.Item(120, 12).Value = "Primera"
$ws.Cells.Item(120, 13).Value = 120
$ws.Cells.Item(120, 14).Value = 12000
$ws.Cells.Item(120, 15).Value = 12500
$ws.Cells.Item(120, 16).Value = 12250
$ws.Cells.Item(120, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(120, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(120, 19).Value = 681
$ws.Cells.Item(120, 20).Value = 18

$ws.Cells.Item(121, 4).Value = 44442
$ws.Cells.Item(121, 12).Value = "Segunda"
$ws.Cells.Item(121, 13).Value = 100
$ws.Cells.Item(121, 14).Value = 11000
$ws.Cells.Item(121, 15).Value = 11500
$ws.Cells.Item(121, 16).Value = 11250
$ws.Cells.Item(121, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(121, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(121, 19).Value = 625
$ws.Cells.Item(121, 20).Value = 18

$ws.Cells.Item(122, 4).Value = 44435
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 560
$ws.Cells.Item(122, 14).Value = 12000
$ws.Cells.Item(122, 15).Value = 12500
$ws.Cells.Item(122, 16).Value = 12250
$ws.Cells.Item(122, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(122, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(122, 19).Value = 681
$ws.Cells.Item(122, 20).Value = 18

$ws.Cells.Item(123, 4).Value = 44435
$ws.Cells.Item(123, 12).Value = "Segunda"
$ws.Cells.Item(123, 13).Value = 400
$ws.Cells.Item(123, 14).Value = 11000
$ws.Cells.Item(123, 15).Value = 11500
$ws.Cells.Item(123, 16).Value = 11250
$ws.Cells.Item(123, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(123, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(123, 19).Value = 625
$ws.Cells.Item(123, 20).Value = 18

$ws.Cells.Item(124, 4).Value = 44319
$ws.Cells.Item(124, 12).Value = "Primera"
$ws.Cells.Item(124, 13).Value = 160
$ws.Cells.Item(124, 14).Value = 10000
$ws.Cells.Item(124, 15).Value = 11000
$ws.Cells.Item(124, 16).Value = 10500
$ws.Cells.Item(124, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(124, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(124, 19).Value = 1050
$ws.Cells.Item(124, 20).Value = 10

$ws.Cells.Item(125, 4).Value = 44319
$ws.Cells.Item(125, 12).Value = "Segunda"
$ws.Cells.Item(125, 13).Value = 120
$ws.Cells.Item(125, 14).Value = 8500
$ws.Cells.Item(125, 15).Value = 9000
$ws.Cells.Item(125, 16).Value = 8750
$ws.Cells.Item(125, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(125, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(125, 19).Value = 875
$ws.Cells.Item(125, 20).Value = 10

$ws.Cells.Item(126, 4).Value = 44376
$ws.Cells.Item(126, 12).Value = "Primera"
$ws.Cells.Item(126, 13).Value = 120
$ws.Cells.Item(126, 14).Value = 10000
$ws.Cells.Item(126, 15).Value = 11000
$ws.Cells.Item(126, 16).Value = 10500
$ws.Cells.Item(126, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(126, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(126, 19).Value = 583
$ws.Cells.Item(126, 20).Value = 18

$ws.Cells.Item(127, 4).Value = 44376
$ws.Cells.Item(127, 12).Value = "Segunda"
$ws.Cells.Item(127, 13).Value = 120
$ws.Cells.Item(127, 14).Value = 8500
$ws.Cells.Item(127, 15).Value = 9000
$ws.Cells.Item(127, 16).Value = 8750
$ws.Cells.Item(127, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(127, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(127, 19).Value = 486
$ws.Cells.Item(127, 20).Value = 18

$ws.Cells.Item(128, 4).Value = 44412
$ws.Cells.Item(128, 12).Value = "Primera"
$ws.Cells.Item(128, 13).Value = 120
$ws.Cells.Item(128, 14).Value = 8000
$ws.Cells.Item(128, 15).Value = 8500
$ws.Cells.Item(128, 16).Value = 8250
$ws.Cells.Item(128, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(128, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(128, 19).Value = 825
$ws.Cells.Item(128, 20).Value = 10

$ws.Cells.Item(129, 4).Value = 44412
$ws.Cells.Item(129, 12).Value = "Segunda"
$ws.Cells.Item(129, 13).Value = 80
$ws.Cells.Item(129, 14).Value = 7000
$ws.Cells.Item(129, 15).Value = 7500
$ws.Cells.Item(129, 16).Value = 7250
$ws.Cells.Item(129, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(129, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(129, 19).Value = 725
$ws.Cells.Item(129, 20).Value = 10

$ws.Cells.Item(130, 4).Value = 44314
$ws.Cells.Item(130, 12).Value = "Primera"
$ws.Cells.Item(130, 13).Value = 120
$ws.Cells.Item(130, 14).Value = 8500
$ws.Cells.Item(130, 15).Value = 9000
$ws.Cells.Item(130, 16).Value = 8750
$ws.Cells.Item(130, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(130, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(130, 19).Value = 875
$ws.Cells.Item(130, 20).Value = 10

$ws.Cells.Item(131, 4).Value = 44399
$ws.Cells.Item(131, 12).Value = "Primera"
$ws.Cells.Item(131, 13).Value = 120
$ws.Cells.Item(131, 14).Value = 10000
$ws.Cells.Item(131, 15).Value = 11000
$ws.Cells.Item(131, 16).Value = 10500
$ws.Cells.Item(131, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(131, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(131, 19).Value = 583
$ws.Cells.Item(131, 20).Value = 18

$ws.Cells.Item(132, 4).Value = 44399
$ws.Cells.Item(132, 12).Value = "Segunda"
$ws.Cells.Item(132, 13).Value = 80
$ws.Cells.Item(132, 14).Value = 9000
$ws.Cells.Item(132, 15).Value = 9000
$ws.Cells.Item(132, 16).Value = 9000
$ws.Cells.Item(132, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(132, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(132, 19).Value = 500
$ws.Cells.Item(132, 20).Value = 18

$ws.Cells.Item(133, 1).Value = 7
$ws.Cells.Item(133, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(133, 3).Value = "Ñuble"
$ws.Cells.Item(133, 4).Value = 44400
$ws.Cells.Item(133, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(133, 5).Value = 16
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100101
$ws.Cells.Item(133, 8).Value = "Berries"
$ws.Cells.Item(133, 9).Value = 100101007
$ws.Cells.Item(133, 10).Value = "Kiwi"
$ws.Cells.Item(133, 11).Value = "Hayward"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 120
$ws.Cells.Item(133, 14).Value = 11000
$ws.Cells.Item(133, 15).Value = 12000
$ws.Cells.Item(133, 16).Value = 11500
$ws.Cells.Item(133, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(133, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(133, 19).Value = 639
$ws.Cells.Item(133, 20).Value = 18

$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value = "Ñuble"
$ws.Cells.Item(134, 4).Value = 44400
$ws.Cells.Item(134, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100101
$ws.Cells.Item(134, 8).Value = "Berries"
$ws.Cells.Item(134, 9).Value = 100101007
$ws.Cells.Item(134, 10).Value = "Kiwi"
$ws.Cells.Item(134, 11).Value = "Hayward"
$ws.Cells.Item(134, 12).Value = "Segunda"
$ws.Cells.Item(134, 13).Value = 120
$ws.Cells.Item(134, 14).Value = 9000
$ws.Cells.Item(134, 15).Value = 10000
$ws.Cells.Item(134, 16).Value = 9500
$ws.Cells.Item(134, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(134, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(134, 19).Value = 528
$ws.Cells.Item(134, 20).Value = 18
